# Auto-generated Excel COM-interop script applying the Unicorn_Profits diff.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) for specific
# rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, matching a scheduled
# market-data refresh run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 638.3077
$ws.Range("I28").Value = 666.5
$ws.Range("J28").Value = 300
$ws.Range("K28").Value = 666.5
$ws.Range("L28").Value = 300
$ws.Range("M28").Value = -181.5
$ws.Range("N28").Value = -1270

# Row 32
$ws.Range("H32").Value = 7744422.5
$ws.Range("J32").Value = 9957079
$ws.Range("L32").Value = 9957079
$ws.Range("N32").Value = -9957731

# Row 53
$ws.Range("H53").Value = 210.2
$ws.Range("I53").Value = 233.45454
$ws.Range("J53").Value = 196.73685
$ws.Range("K53").Value = 233.45454
$ws.Range("L53").Value = 196.73685
$ws.Range("M53").Value = 403.54546
$ws.Range("N53").Value = -1470.73685

# Row 88
$ws.Range("H88").Value = 3693.9048
$ws.Range("I88").Value = 1124.5
$ws.Range("J88").Value = 4298.4707
$ws.Range("K88").Value = 1124.5
$ws.Range("L88").Value = 4298.4707
$ws.Range("M88").Value = -718.5
$ws.Range("N88").Value = -5110.4707

# Row 91
$ws.Range("H91").Value = 3693.9048
$ws.Range("I91").Value = 1124.5
$ws.Range("J91").Value = 4298.4707
$ws.Range("K91").Value = 1124.5
$ws.Range("L91").Value = 4298.4707
$ws.Range("M91").Value = 279.5
$ws.Range("N91").Value = -7106.4707

# Row 96
$ws.Range("H96").Value = 1284.4
$ws.Range("I96").Value = 888.93335
$ws.Range("K96").Value = 2666.80005
$ws.Range("M96").Value = -1293.80005

# Row 132
$ws.Range("H132").Value = 3060.1025
$ws.Range("I132").Value = 1946.3793
$ws.Range("J132").Value = 6289.9
$ws.Range("K132").Value = 5839.1379
$ws.Range("L132").Value = 18869.7
$ws.Range("M132").Value = -3309.1379
$ws.Range("N132").Value = -23929.7

# Row 137
$ws.Range("H137").Value = 2117.3606
$ws.Range("I137").Value = 1972.3265
$ws.Range("J137").Value = 2709.5833
$ws.Range("K137").Value = 5916.979499999999
$ws.Range("L137").Value = 8128.749899999999
$ws.Range("M137").Value = -3366.979499999999
$ws.Range("N137").Value = -13228.7499

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1112.6578
$ws.Range("I45").Value = 923.7273
$ws.Range("K45").Value = 923.7273
$ws.Range("M45").Value = -546.7273

# Row 122
$ws.Range("H122").Value = 6887.125
$ws.Range("I122").Value = 7024.5
$ws.Range("J122").Value = 6475
$ws.Range("K122").Value = 21073.5
$ws.Range("L122").Value = 19425
$ws.Range("M122").Value = -18623.5
$ws.Range("N122").Value = -24325

# Row 132
$ws.Range("H132").Value = 2770.9656
$ws.Range("I132").Value = 3349.6
$ws.Range("J132").Value = 2151
$ws.Range("K132").Value = 10048.8
$ws.Range("L132").Value = 6453
$ws.Range("M132").Value = -7518.799999999999
$ws.Range("N132").Value = -11513

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 4600.75
$ws.Range("I86").Value = 6551.2
$ws.Range("J86").Value = 1350
$ws.Range("K86").Value = 6551.2
$ws.Range("L86").Value = 1350
$ws.Range("M86").Value = -5428.2
$ws.Range("N86").Value = -3596

# Row 89
$ws.Range("H89").Value = 4600.75
$ws.Range("I89").Value = 6551.2
$ws.Range("J89").Value = 1350
$ws.Range("K89").Value = 32756
$ws.Range("L89").Value = 6750
$ws.Range("M89").Value = -27140
$ws.Range("N89").Value = -17982

# Row 99
$ws.Range("H99").Value = 5774604.5
$ws.Range("I99").Value = 1837476.5
$ws.Range("J99").Value = 33334500
$ws.Range("K99").Value = 1837476.5
$ws.Range("L99").Value = 33334500
$ws.Range("M99").Value = -1835978.5
$ws.Range("N99").Value = -33337496

# Row 107
$ws.Range("H107").Value = 2057.6875
$ws.Range("I107").Value = 1510
$ws.Range("J107").Value = 4431
$ws.Range("K107").Value = 1510
$ws.Range("L107").Value = 4431
$ws.Range("M107").Value = 410
$ws.Range("N107").Value = -8271

$ws = $wb.Worksheets.Item("CRP")
# Row 132
$ws.Range("H132").Value = 1518.3016
$ws.Range("I132").Value = 999.8200000000001
$ws.Range("J132").Value = 3512.4614
$ws.Range("K132").Value = 2999.46
$ws.Range("L132").Value = 10537.3842
$ws.Range("M132").Value = -469.46
$ws.Range("N132").Value = -15597.3842

$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Range("H14").Value = 87
$ws.Range("I14").Value = 87
$ws.Range("K14").Value = 261
$ws.Range("M14").Value = -88

# Row 38
$ws.Range("H38").Value = 17543926
$ws.Range("I38").Value = 63.75
$ws.Range("J38").Value = 30303098
$ws.Range("K38").Value = 191.25
$ws.Range("L38").Value = 90909294
$ws.Range("M38").Value = 155.75
$ws.Range("N38").Value = -90909988

# Row 116
$ws.Range("H116").Value = 3539.2083
$ws.Range("I116").Value = 581.75
$ws.Range("J116").Value = 6496.6665
$ws.Range("K116").Value = 1745.25
$ws.Range("L116").Value = 19489.9995
$ws.Range("M116").Value = 1696.75
$ws.Range("N116").Value = -26373.9995

# Row 131
$ws.Range("H131").Value = 1951.6986
$ws.Range("I131").Value = 2490.8948
$ws.Range("J131").Value = 1761.9814
$ws.Range("K131").Value = 7472.6844
$ws.Range("L131").Value = 5285.9442
$ws.Range("M131").Value = -2432.6844
$ws.Range("N131").Value = -15365.9442

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 5263
$ws.Range("I80").Value = 6243.8887
$ws.Range("J80").Value = 3791.6667
$ws.Range("K80").Value = 6243.8887
$ws.Range("L80").Value = 3791.6667
$ws.Range("M80").Value = -5245.8887
$ws.Range("N80").Value = -5787.6667

# Row 83
$ws.Range("H83").Value = 5263
$ws.Range("I83").Value = 6243.8887
$ws.Range("J83").Value = 3791.6667
$ws.Range("K83").Value = 31219.4435
$ws.Range("L83").Value = 18958.3335
$ws.Range("M83").Value = -26227.4435
$ws.Range("N83").Value = -28942.3335

# Row 122
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2226.2222
$ws.Range("I7").Value = 2000.9286
$ws.Range("J7").Value = 2468.8462
$ws.Range("K7").Value = 2000.9286
$ws.Range("L7").Value = 2580.7144
$ws.Range("M7").Value = -1888.9286
$ws.Range("N7").Value = -2692.8462

# Row 55
$ws.Range("H55").Value = 542.9
$ws.Range("I55").Value = 288.16666
$ws.Range("J55").Value = 925
$ws.Range("K55").Value = 288.16666
$ws.Range("L55").Value = 925
$ws.Range("M55").Value = -115.16666
$ws.Range("N55").Value = -1271

# Row 82
$ws.Range("H82").Value = 2645.2727
$ws.Range("I82").Value = 2591.3333
$ws.Range("K82").Value = 2591.3333
$ws.Range("M82").Value = -2230.3333

# Row 85
$ws.Range("H85").Value = 2645.2727
$ws.Range("I85").Value = 2591.3333
$ws.Range("K85").Value = 2591.3333
$ws.Range("M85").Value = -1343.3333

# Row 126
$ws.Range("H126").Value = 2226.2222
$ws.Range("I126").Value = 2000.9286
$ws.Range("J126").Value = 2468.8462
$ws.Range("K126").Value = 6002.7858
$ws.Range("L126").Value = 7406.5386
$ws.Range("M126").Value = -3532.7858
$ws.Range("N126").Value = -12346.5386

# Row 132
$ws.Range("H132").Value = 5696.0166
$ws.Range("I132").Value = 2109.1282
$ws.Range("J132").Value = 12357.381
$ws.Range("K132").Value = 6327.3846
$ws.Range("L132").Value = 37072.143
$ws.Range("M132").Value = -3797.3846
$ws.Range("N132").Value = -42132.143

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1858.9744
$ws.Range("I132").Value = 1122
$ws.Range("J132").Value = 3175
$ws.Range("K132").Value = 3366
$ws.Range("L132").Value = 9525
$ws.Range("M132").Value = -836
$ws.Range("N132").Value = -14585

